$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.414.47"
$ws.Range("D3").Value = "2.456.76"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.14"
$ws.Range("E5").Value = "  -2.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.90"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "2.455.94"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("E10").Value = "  -6.85%  "
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("E12").Value = "  -6.26%  "
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Value = "2.906.54"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "68.330.78"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.33"
$ws.Range("E17").Value = "  -5.66%  "
$ws.Range("D18").Value = "2.479.96"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.93"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.20"
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.15"
$ws.Range("E21").Value = "  -4.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.76"
$ws.Range("E22").Value = "  -3.78%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -3.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.41"
$ws.Range("E25").Value = "  -4.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.09"
$ws.Range("E26").Value = "  +7.87%  "
$ws.Range("E27").Value = "  -6.23%  "
$ws.Range("E28").Value = "  -1.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.10"
$ws.Range("E29").Value = "  -7.33%  "
$ws.Range("D30").Value = "0.0₃0828"
$ws.Range("E30").Value = "  -6.98%  "
$ws.Range("E31").Value = "  -8.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.44"
$ws.Range("E32").Value = "  +133.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "430.24"
$ws.Range("E34").Value = "  -6.07%  "
$ws.Range("E35").Value = "  -3.67%  "
$ws.Range("E36").Value = "  -4.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.43"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.01"
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.110"
$ws.Range("E40").Value = "  -5.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.86"
$ws.Range("E41").Value = "  -2.68%  "
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("E43").Value = "  -5.38%  "
$ws.Range("E44").Value = "  -5.86%  "
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  -6.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "134.61"
$ws.Range("E47").Value = "  -4.74%  "
$ws.Range("E48").Value = "  -4.10%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.482"
$ws.Range("E50").Value = "  -7.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.564"
$ws.Range("E51").Value = "  -2.63%  "
